$d = $word.ActiveDocument

# 1. Merge the whitespace runs before "End Time: " into a single run
#    ("   " + "                    " + "End Time: " -> 23 spaces + "End Time: ")
$d.Content.Find.Execute("   " + "                    " + "End Time: ", $false, $false, $false, $false, $false, $true, 1, $false, "                       End Time: ", 2)

# 2. Renumber the US_ task ranges in "Tasks Assigned:"
#    US_17 - US_19 -> US_17 - US_18
$d.Content.Find.Execute("US_17 – US_19", $false, $false, $false, $false, $false, $true, 1, $false, "US_17 – US_18", 2)

#    US_20 - US_22 -> US_19 - US_21
$d.Content.Find.Execute("US_20 – US_22", $false, $false, $false, $false, $false, $true, 1, $false, "US_19 – US_21", 2)

#    US_23 - US_25 -> US_22 - US_23
$d.Content.Find.Execute("US_23 – US_25", $false, $false, $false, $false, $false, $true, 1, $false, "US_22 – US_23", 2)
